$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Gender values for rows 3 and 4 (swap male/female)
$ws.Range("E3").Value = "male"
$ws.Range("E4").Value = "female"

# Update the active selection to F10
$ws.Range("F10").Select()
